# Apply the edits described by the diff:
#  - RegistrationForm!B6 changes from a date serial value to the text "1998,Jan,1"
#  - RegistrationForm sheet becomes the active/selected sheet, with B7 selected
#  - Login sheet is no longer the selected sheet

$wb = $excel.ActiveWorkbook

$wsRegistration = $wb.Worksheets.Item("RegistrationForm")

# Replace the DOB cell value (was a date serial number) with the literal text value
$wsRegistration.Range("B6").Value = "1998,Jan,1"

# Make RegistrationForm the active sheet and select cell B7 on it,
# so the sheet view / selection metadata matches the target workbook.
$wsRegistration.Activate() | Out-Null
$wsRegistration.Range("B7").Select() | Out-Null
